$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")

# Fill in row 36 (previously only had B36 = 6977)
$ws.Range("C36").Value = 43925
$ws.Range("D36").Value = 0.94791666666666663
$ws.Range("E36").Value = 0.96180555555555547
$ws.Range("G36").Value = "Created RTL netlist images and Post-fit images"

# Fill in row 37 (previously empty)
$ws.Range("B37").Value = 6977
$ws.Range("C37").Value = 43925
$ws.Range("D37").Value = 0.96180555555555547
$ws.Range("E37").Value = 0.97916666666666663
$ws.Range("G37").Value = "Updated Timing waveforms for LogicUnit.vhd"

# Update G35 text (was "Updated Functional Waveforms.") to the new, more specific text
$ws.Range("G35").Value = "Updated Functional Waveforms for LogicUnit.vhd"

# Update selection to F37
$ws.Range("F37").Select()
